$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original styles for columns D and E (data rows 2-51) so that
# forcing a text NumberFormat does not alter their (default) cell style.
$origStyleD = $ws.Range("D2:D51").Style
$origStyleE = $ws.Range("E2:E51").Style

# Force column D and E (rows 2-51) to text format so that numeric-looking
# strings (e.g. "8.75", "0.999") are NOT auto-converted to numbers by Excel,
# matching the inline/shared string cell type used in the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "62.652.50"
$ws.Range("E2").Value = "  -7.33%  "

# Row 3
$ws.Range("D3").Value = "3.204.90"
$ws.Range("E3").Value = "  -8.83%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "508.61"
$ws.Range("E5").Value = "  -8.00%  "

# Row 6
$ws.Range("D6").Value = "167.45"
$ws.Range("E6").Value = "  -15.33%  "

# Row 7
$ws.Range("D7").Value = "0.584"
$ws.Range("E7").Value = "  -8.27%  "

# Row 8
$ws.Range("E8").Value = "  +0.15%  "

# Row 9
$ws.Range("D9").Value = "3.202.58"
$ws.Range("E9").Value = "  -8.69%  "

# Row 10
$ws.Range("D10").Value = "0.587"
$ws.Range("E10").Value = "  -10.44%  "

# Row 11
$ws.Range("D11").Value = "53.71"
$ws.Range("E11").Value = "  -11.87%  "

# Row 12
$ws.Range("E12").Value = "  -9.84%  "

# Row 13
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  -7.72%  "

# Row 14
$ws.Range("D14").Value = "8.75"
$ws.Range("E14").Value = "  -10.71%  "

# Row 15
$ws.Range("D15").Value = "3.734.73"
$ws.Range("E15").Value = "  -8.32%  "

# Row 16
$ws.Range("D16").Value = "3.214.78"
$ws.Range("E16").Value = "  -8.45%  "

# Row 17
$ws.Range("D17").Value = "62.630.72"
$ws.Range("E17").Value = "  -6.95%  "

# Row 18
$ws.Range("D18").Value = "0.112"
$ws.Range("E18").Value = "  -9.45%  "

# Row 19
$ws.Range("D19").Value = "16.82"
$ws.Range("E19").Value = "  -8.26%  "

# Row 20
$ws.Range("D20").Value = "10.62"
$ws.Range("E20").Value = "  -10.09%  "

# Row 21
$ws.Range("D21").Value = "0.927"
$ws.Range("E21").Value = "  -9.57%  "

# Row 22
$ws.Range("D22").Value = "359.47"
$ws.Range("E22").Value = "  -8.70%  "

# Row 23
$ws.Range("D23").Value = "3.63"
$ws.Range("E23").Value = "  -8.54%  "

# Row 24
$ws.Range("D24").Value = "78.08"
$ws.Range("E24").Value = "  -8.41%  "

# Row 25
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "10.55"
$ws.Range("E25").Value = "  -11.13%  "

# Row 26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "6.05"
$ws.Range("E26").Value = "  -1.70%  "

# Row 27
$ws.Range("D27").Value = "3.79"
$ws.Range("E27").Value = "  -1.87%  "

# Row 28
$ws.Range("D28").Value = "2.58"
$ws.Range("E28").Value = "  -8.25%  "

# Row 29
$ws.Range("D29").Value = "10.91"
$ws.Range("E29").Value = "  -11.42%  "

# Row 30
$ws.Range("D30").Value = "8.02"
$ws.Range("E30").Value = "  -9.44%  "

# Row 31
$ws.Range("D31").Value = "27.76"
$ws.Range("E31").Value = "  -11.54%  "

# Row 32
$ws.Range("D32").Value = "612.95"
$ws.Range("E32").Value = "  -14.88%  "

# Row 33
$ws.Range("D33").Value = "6.33"
$ws.Range("E33").Value = "  -10.07%  "

# Row 34
$ws.Range("D34").Value = "10.89"
$ws.Range("E34").Value = "  -6.77%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.101"
$ws.Range("E35").Value = "  -8.51%  "

# Row 36
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "56.84"
$ws.Range("E36").Value = "  -11.16%  "

# Row 37
$ws.Range("E37").Value = "  -0.10%  "

# Row 38
$ws.Range("D38").Value = "35.58"
$ws.Range("E38").Value = "  -7.35%  "

# Row 39
$ws.Range("D39").Value = "0.370"
$ws.Range("E39").Value = "  -5.56%  "

# Row 40
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.16%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0659"
$ws.Range("E41").Value = "  -3.09%  "

# Row 42
$ws.Range("D42").Value = "2.833.96"
$ws.Range("E42").Value = "  -7.51%  "

# Row 43
$ws.Range("D43").Value = "0.118"
$ws.Range("E43").Value = "  -9.65%  "

# Row 44
$ws.Range("D44").Value = "2.36"
$ws.Range("E44").Value = "  -6.06%  "

# Row 45
$ws.Range("D45").Value = "2.57"
$ws.Range("E45").Value = "  -7.36%  "

# Row 46
$ws.Range("D46").Value = "2.56"
$ws.Range("E46").Value = "  -14.91%  "

# Row 47
$ws.Range("D47").Value = "0.0378"
$ws.Range("E47").Value = "  -7.12%  "

# Row 48
$ws.Range("D48").Value = "2.93"
$ws.Range("E48").Value = "  -1.26%  "

# Row 49
$ws.Range("E49").Value = "  +2.22%  "

# Row 50
$ws.Range("D50").Value = "132.91"
$ws.Range("E50").Value = "  -4.21%  "

# Row 51
$ws.Range("E51").Value = "  -8.05%  "

# Restore original styles
$ws.Range("D2:D51").Style = $origStyleD
$ws.Range("E2:E51").Style = $origStyleE
